$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value  = 107
$ws.Range("I2").Value  = 251
$ws.Range("J2").Value  = 1099
$ws.Range("K2").Value  = 4
$ws.Range("L2").Value  = 314
$ws.Range("M2").Value  = 18
$ws.Range("N2").Value  = 202
$ws.Range("O2").Value  = 0
$ws.Range("P2").Value  = 3
$ws.Range("Q2").Value  = 2
$ws.Range("R2").Value  = 13
$ws.Range("S2").Value  = 108
$ws.Range("T2").Value  = 212
$ws.Range("U2").Value  = 23
$ws.Range("V2").Value  = 1659
$ws.Range("W2").Value  = 1
$ws.Range("X2").Value  = 1743
$ws.Range("Y2").Value  = 2
$ws.Range("Z2").Value  = 23
$ws.Range("AA2").Value = 10

$wb.Save()
